$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column widths: C (3), D (4), H (8)
# (ColumnWidth is offset from the stored XML "width" by ~5/6 of a character
#  for this font, so subtract that to land exactly on the target stored widths)
$ws.Columns.Item(3).ColumnWidth = 45.166666666666664
$ws.Columns.Item(4).ColumnWidth = 55.166666666666664
$ws.Columns.Item(8).ColumnWidth = 55.166666666666664

# Full target data for rows 2-9 (columns A-H)
$data = @(
    @("1330823", "https://aiesec.org/opportunity/global-talent/1330823", "Accelerate Romania | Marketing & SEO Intern", "Timișoara, Romania", "No", "1 applicant", "9 - 12 Weeks", "Maschinenbau"),
    @("1330822", "https://aiesec.org/opportunity/global-talent/1330822", "Accelerate Romania | Web Development Intern", "Timișoara, Romania", "No", "1 applicant", "9 - 12 Weeks", "Maschinenbau"),
    @("1330821", "https://aiesec.org/opportunity/global-talent/1330821", "Foreign Trade Expert", "Esenyurt, Balıkyolu, 34510 Esenyurt/İstanbul, Türkiye", "No", "5 applicants", "6 - 18 Months", "Yelken Kalıp Pencere &Kapı AKS. ve MET. SAN. TİC. AŞ."),
    @("1330755", "https://aiesec.org/opportunity/global-talent/1330755", "Export Specialist", "Konya, Türkiye", "No", "1 applicant", "6 - 18 Months", "ŞÖLEN MEDİKAL GIDA TEMİZLİK İNŞAAT LİMİTED ŞİRKETİ"),
    @("1330716", "https://aiesec.org/opportunity/global-talent/1330716", "Sales and Marketing Intern", "Boralesgamuwa, Sri Lanka", "No", "0 applicants", "9 - 12 Weeks", "No Name Collective Social Action"),
    @("1330709", "https://aiesec.org/opportunity/global-talent/1330709", "Export Specialist", "Konya, Türkiye", "No", "5 applicants", "6 - 18 Months", "Anıl Yatağanlı A.Ş"),
    @("1330185", "https://aiesec.org/opportunity/global-talent/1330185", "Front-End Web Developer", "Ciudad Juárez, Chihuahua, Mexico", "No", "20 applicants", "6 - 18 Months", "EP&O Corporation"),
    @("1329565", "https://aiesec.org/opportunity/global-talent/1329565", "Marketing", "Sfax, Tunisia", "No", "6 applicants", "9 - 12 Weeks", "SRTC Country tunisia city sfax")
)

# Force column A (opportunity ID) to be stored as text, matching the source data
$ws.Range("A2:A9").NumberFormat = "@"

$rowIndex = 2
foreach ($row in $data) {
    $colIndex = 1
    foreach ($value in $row) {
        $ws.Cells.Item($rowIndex, $colIndex).Value = $value
        $colIndex++
    }
    $rowIndex++
}
